# Insert a new data row at row 58 (shifting existing rows 58:91 down to 59:92)
# and populate it with the new "Arveja Verde" observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(58).Insert()

$ws.Cells.Item(58, 1).Value = 4
$ws.Cells.Item(58, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(58, 3).Value = "Los Lagos"
$ws.Cells.Item(58, 4).Value = 44529
$ws.Cells.Item(58, 5).Value = 10
$ws.Cells.Item(58, 6).Value = 100112022
$ws.Cells.Item(58, 7).Value = "Arveja Verde"
$ws.Cells.Item(58, 8).Value = "Sin especificar"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 35
$ws.Cells.Item(58, 11).Value = 17000
$ws.Cells.Item(58, 12).Value = 17000
$ws.Cells.Item(58, 13).Value = 17000
$ws.Cells.Item(58, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(58, 15).Value = "Región del Maule"
$ws.Cells.Item(58, 16).Value = 680
$ws.Cells.Item(58, 17).Value = 25
$ws.Cells.Item(58, 18).Value = "Hortaliza"
